$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 1: new label in A1 ("Total horas: 75.5")
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Total horas: 75.5"

# ---------------------------------------------------------------
# Row 2: extend the date header from G2 out to W2
# ---------------------------------------------------------------
$ws.Range("G2").Copy()
$ws.Range("H2:W2").PasteSpecial(-4122)

$ws.Range("H2").Value = Get-Date -Year 2020 -Month 6 -Day 11 -Hour 0 -Minute 0 -Second 0
$ws.Range("I2").Value = Get-Date -Year 2020 -Month 6 -Day 14 -Hour 0 -Minute 0 -Second 0
$ws.Range("J2").Value = Get-Date -Year 2020 -Month 6 -Day 17 -Hour 0 -Minute 0 -Second 0
$ws.Range("K2").Value = Get-Date -Year 2020 -Month 6 -Day 18 -Hour 0 -Minute 0 -Second 0
$ws.Range("L2").Value = Get-Date -Year 2020 -Month 6 -Day 19 -Hour 0 -Minute 0 -Second 0
$ws.Range("M2").Value = Get-Date -Year 2020 -Month 6 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Range("N2").Value = Get-Date -Year 2020 -Month 6 -Day 21 -Hour 0 -Minute 0 -Second 0
$ws.Range("O2").Value = Get-Date -Year 2020 -Month 6 -Day 23 -Hour 0 -Minute 0 -Second 0
$ws.Range("P2").Value = Get-Date -Year 2020 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0
$ws.Range("Q2").Value = Get-Date -Year 2020 -Month 6 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Range("R2").Value = Get-Date -Year 2020 -Month 6 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("S2").Value = Get-Date -Year 2020 -Month 6 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Range("T2").Value = Get-Date -Year 2020 -Month 6 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("U2").Value = Get-Date -Year 2020 -Month 6 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("V2").Value = Get-Date -Year 2020 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("W2").Value = Get-Date -Year 2020 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0

# ---------------------------------------------------------------
# Row 9: rename the task from Deep Q-Learning to Q-Learning
# ---------------------------------------------------------------
$ws.Range("A9").Value = "Aprendizaje por refuerzo (Q-Learning)"

# ---------------------------------------------------------------
# Insert a fresh row at 10 (old row 10 "Implementación..." becomes
# row 11, old row 11 "Estructuras básicas..." becomes row 12)
# ---------------------------------------------------------------
$ws.Rows("10:10").Insert()
# the insert drags down D:F formatting from row 9 into the blank new
# row - strip that back out before laying down the real content
$ws.Range("D10:F10").Clear()

# Row 10: "Aprendizaje por refuerzo (Actor-Critic)" sub task
$ws.Range("A10").Value = "Aprendizaje por refuerzo (Actor-Critic)"

$ws.Range("F9").Copy()
$ws.Range("L10:N10").PasteSpecial(-4122)
$ws.Range("L10").Value = "1 h."
$ws.Range("M10").Value = "2 h."
$ws.Range("N10").Value = "3 h."

# Row 11 (old row 10): section header renamed to "Implementación Q-Learning"
$ws.Range("A11").Value = "Implementación Q-Learning"

# Row 12 (old row 11): "Estructuras básicas (estados y agente)" - unchanged,
# only needs 3 new blank-but-styled cells J12:L12 added
$ws.Range("J12").Interior.Color = $ws.Range("F9").Interior.Color
$ws.Range("K12").Interior.Color = $ws.Range("F9").Interior.Color
$ws.Range("L12").Interior.Color = $ws.Range("F9").Interior.Color

# ---------------------------------------------------------------
# Insert a fresh row at 13 (new "Algoritmo Q-learning" row)
# ---------------------------------------------------------------
$ws.Rows("13:13").Insert()
# the insert drags down G formatting from row 12 into the blank new row
$ws.Range("G13").Clear()

$ws.Range("A13").Value = "Algoritmo Q-learning"

$ws.Range("F9").Copy()
$ws.Range("H13:L13").PasteSpecial(-4122)
$ws.Range("H13").Value = "3.5 h."
$ws.Range("I13").Value = "4.5 h."
$ws.Range("J13").Value = "4.5 h."
$ws.Range("K13").Value = "2.5 h."
$ws.Range("L13").Value = "2.5 h."

# ---------------------------------------------------------------
# New rows 14-17: "Servidor" section
# ---------------------------------------------------------------
$ws.Range("A8").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Servicios en línea"

$ws.Range("A9").Copy()
$ws.Range("A15:A17").PasteSpecial(-4122)
$ws.Range("A15").Value = "Estructuras básicas (conexiones)"
$ws.Range("A16").Value = "Implementación partida en red"
$ws.Range("A17").Value = "Servidor"

$ws.Range("F9").Copy()
$ws.Range("O15:U15").PasteSpecial(-4122)
$ws.Range("O15").Value = "3 h."
$ws.Range("P15").Value = "4 h."
$ws.Range("S15").Value = "4 h."
$ws.Range("T15").Value = "3 h."
$ws.Range("U15").Value = "2 h."

$ws.Range("F9").Copy()
$ws.Range("Q16:R16").PasteSpecial(-4122)
$ws.Range("Q16").Value = "5 h."
$ws.Range("R16").Value = "4 h."

$ws.Range("F9").Copy()
$ws.Range("U17:W17").PasteSpecial(-4122)
$ws.Range("U17").Value = "2 h."
$ws.Range("V17").Value = "3.5 h."
$ws.Range("W17").Value = "4 h."

$ws.Range("T27").Select()
